# Jahresplan 2024/25 IM23d Meyer - "Add files via upload" edit
#
# Content changes applied:
#  - Masterplan!E28 / Masterplan!F28 were blank; the author filled them in
#    (E28 = "?" as a placeholder note, F28 = "PPL" i.e. "Physisches
#    Präsenzlernen", matching the Lernformen picklist already used in the
#    surrounding rows of this block, e.g. F25/F27).
#  - The active selection on the Masterplan sheet moved from N23 to M27
#    (reflecting where the author was last working/clicking).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Masterplan")

$ws.Range("E28").Value = "?"
$ws.Range("F28").Value = "PPL"

$ws.Activate() | Out-Null
$ws.Range("M27").Select() | Out-Null
